{"js": "// Append three empty paragraphs followed by a new paragraph of text to the\n// end of the document body (right before the final section break), matching\n// the \"grid search visuals done\" commit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the current last paragraph (\"...However the code is easily\n// adaptable\") and insert the new content after it, in order.\nlet anchor = paragraphs.getLast();\n\nconst blank1 = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst blank2 = blank1.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst blank3 = blank2.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nblank3.insertParagraph(\n  \"Test how top k topics picked for the bipartite network construction affects matrix similarity at the end. An check if more k aka more recombinations are cool\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Append three empty paragraphs followed by a new paragraph of text to the\n# end of the document body (right before the final section break), matching\n# the \"grid search visuals done\" commit.\n\n$d = $word.ActiveDocument\n\n# Three blank paragraphs after the current last paragraph (\"...However the\n# code is easily adaptable\"), plus one more paragraph break to hold the new\n# note text (so it lands in its own, fourth, new paragraph).\n$p = $d.Paragraphs.Last\n$p.Range.InsertParagraphAfter()\n\n$p = $d.Paragraphs.Last\n$p.Range.InsertParagraphAfter()\n\n$p = $d.Paragraphs.Last\n$p.Range.InsertParagraphAfter()\n\n$p = $d.Paragraphs.Last\n$p.Range.InsertParagraphAfter()\n\n# Fourth new (now-last) paragraph gets the new note text.\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"Test how top k topics picked for the bipartite network construction affects matrix similarity at the end. An check if more k aka more recombinations are cool\")\n"}
